# Applies the "Automatic update of files" edit:
#  1. Updates the "Förändrad" (C column) date for every data row (2-157)
#     from 45184 to 45186.
#  2. Adds a friendly display-text second argument (the "Beteckning" in
#     column A) to the HYPERLINK() formulas found in columns S, T, V, W,
#     X, Y for the rows that have them (rows 2-6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 157
$newDate = 45186

# 1) Update the "Förändrad" date column (C) for every data row.
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $ws.Cells.Item($r, 3).Value = $newDate
}

# 2) Add the Beteckning text as the second HYPERLINK() argument for the
#    link columns S (19), T (20), V (22), W (23), X (24), Y (25).
$linkCols = @(19, 20, 22, 23, 24, 25)

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $name = $ws.Cells.Item($r, 1).Value2
    if (-not $name) { continue }

    foreach ($c in $linkCols) {
        $cell = $ws.Cells.Item($r, $c)
        $f = $cell.Formula
        if ($f -and $f.Length -gt 0 -and $f.ToUpper().Contains("HYPERLINK") -and -not $f.Contains(",")) {
            $newFormula = $f.Substring(0, $f.Length - 1) + ', "' + $name + '")'
            $cell.Formula = $newFormula
        }
    }
}
